$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 29

# Column A holds a date string like the other rows ("MM/DD/YYYY" stored as
# text, not a real date). Force text formatting first so Excel's
# autodetection doesn't silently convert it into a date serial number,
# then drop back to the default style so the new row's cells don't pick up
# a stray custom format.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "12/23/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 12280.2
$ws.Cells.Item($row, 3).Value = 0.2057852427978976
$ws.Cells.Item($row, 4).Value = 0.7942147572021024
$ws.Cells.Item($row, 5).Value = -131.61
$ws.Cells.Item($row, 6).Value = -27.14
$ws.Cells.Item($row, 7).Value = -20788.65
$ws.Cells.Item($row, 8).Value = -68.06999999999999
$ws.Cells.Item($row, 9).Value = -475.77
$ws.Cells.Item($row, 10).Value = -15.84
